$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "3200 Progress Ln Eureka, CA 92000"
$ws.Range("D3").Value = "(321) 555-4321"
$ws.Range("E4").Value = "3200 Progress Ln Eureka, CA 92000"
$ws.Range("D5").Value = "(321) 555-4321"
